$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8650946021080017
$ws.Range("B1").Value = 1.15207040309906
$ws.Range("C1").Value = 5.391704559326172
$ws.Range("D1").Value = 1.492438912391663
$ws.Range("E1").Value = 0.8688130974769592
